$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.983.46'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.144.81'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.61'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.61'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.140.88'
$ws.Range("E8").Value = '  +2.31%  '
$ws.Range("E9").Value = '  +4.21%  '
$ws.Range("E10").Value = '  +6.43%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.503'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +7.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +12.30%  '
$ws.Range("E14").Value = '  +6.91%  '
$ws.Range("D15").Value = '3.658.04'
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").Value = '65.000.14'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("E17").Value = '  +6.48%  '
$ws.Range("D18").Value = '3.150.28'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '509.98'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +6.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.88'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +7.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.729'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +8.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.50'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +13.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.84'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.47'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.18%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +4.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.73'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +8.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.04%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  +4.35%  '
$ws.Range("E33").Value = '  +6.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.02'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +8.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.57'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.59'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '473.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.54%  '
$ws.Range("E38").Value = '  +4.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0856'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.03'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("D41").Value = '3.107.83'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.60'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("E43").Value = '  +4.14%  '
$ws.Range("E44").Value = '  +11.91%  '
$ws.Range("E45").Value = '  +13.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.13'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.03%  '
$ws.Range("D47").Value = '0.0₃0579'
$ws.Range("E47").Value = '  +12.52%  '
$ws.Range("E49").Value = '  +3.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +10.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.65'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.96%  '
